# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (data row): date serial number
$ws.Range("A2").Value = 46021

# Hourly prices B2:Z2
$ws.Range("B2").Value = 104.79
$ws.Range("C2").Value = 97.67
$ws.Range("D2").Value = 87.34
$ws.Range("E2").Value = 85.48
$ws.Range("F2").Value = 82.09
$ws.Range("G2").Value = 82.36
$ws.Range("H2").Value = 84.2
$ws.Range("I2").Value = 90.11
$ws.Range("J2").Value = 99.61
$ws.Range("K2").Value = 94.97
$ws.Range("L2").Value = 85.92
$ws.Range("M2").Value = 79.38
$ws.Range("N2").Value = 76.16
$ws.Range("O2").Value = 73.04000000000001
$ws.Range("P2").Value = 73.03
$ws.Range("Q2").Value = 75.84999999999999
$ws.Range("R2").Value = 87.15000000000001
$ws.Range("S2").Value = 99.59999999999999
$ws.Range("T2").Value = 104.9
$ws.Range("U2").Value = 105.77
$ws.Range("V2").Value = 105.27
$ws.Range("W2").Value = 106.52
$ws.Range("X2").Value = 104.75
$ws.Range("Y2").Value = 103.71
$ws.Range("Z2").Value = 91.23999999999999

# Slot max/price columns
$ws.Range("AB2").Value = 105.06
$ws.Range("AD2").Value = 105.9
$ws.Range("AF2").Value = 105.34

# Slot_min_price label
$ws.Range("AG2").Value = "2h-16h"
